# Update crypto price/volume snapshot data (scraped refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain decimals (e.g. "706.87") which Excel
# would otherwise auto-convert to a number on assignment; force those cells
# to Text format first so the value is stored as a literal string, matching
# the existing column formatting (all Price/Volume cells are text).
$ws.Range("D2").Value = "71.182.27"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.813.37"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "706.87"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.86"
$ws.Range("D7").Value = "3.811.95"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.95"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "4.457.16"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "3.802.29"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "71.169.37"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.17"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.45"
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.19"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "3.965.63"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.41"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.41"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.12"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.17"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "3.775.46"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.27"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "170.65"
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.61"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "422.83"
$ws.Range("E49").Value = "  +4.71%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  -1.50%  "
